$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear cells that no longer hold values
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Update the remaining values for rows 7-19
$ws.Range("C7").Value = 1.097054137926201
$ws.Range("E7").Value = 1.37755776875883

$ws.Range("C8").Value = 1.385527545913412
$ws.Range("E8").Value = 1.329814931661888

$ws.Range("C9").Value = 1.296301936385214
$ws.Range("E9").Value = 1.355477993452414

$ws.Range("C10").Value = 2.441628883342295
$ws.Range("E10").Value = 1.386547975635688

$ws.Range("C11").Value = 2.565764046666463
$ws.Range("E11").Value = 1.833587970352424

$ws.Range("C12").Value = 1.263447557103259
$ws.Range("E12").Value = 1.485511920344451

$ws.Range("C13").Value = 2.117022522597423
$ws.Range("E13").Value = 1.745834498329324

$ws.Range("C14").Value = 2.149400276001101
$ws.Range("E14").Value = 1.76475225558832

$ws.Range("C15").Value = 2.453568910971748
$ws.Range("E15").Value = 2.131436976903012

$ws.Range("C16").Value = 0.812682184439506
$ws.Range("E16").Value = 1.556352278772266

$ws.Range("C17").Value = 0.9940067218177528
$ws.Range("E17").Value = 1.820779918499094

$ws.Range("C18").Value = 1.634555928116921
$ws.Range("E18").Value = 1.554016159863814

$ws.Range("C19").Value = 0.6231570351797
$ws.Range("E19").Value = 1.581524829939718
